$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (data as it existed before the edit)
# Row 1 is header, row 20 is unchanged. Rows 2-19,21-22 are permuted.
$mapping = @{
    2  = 8
    3  = 10
    4  = 4
    5  = 19
    6  = 2
    7  = 16
    8  = 11
    9  = 12
    10 = 13
    11 = 14
    12 = 21
    13 = 22
    14 = 15
    15 = 3
    16 = 5
    17 = 6
    18 = 7
    19 = 9
    20 = 20
    21 = 17
    22 = 18
}

# Columns whose values participate in the row permutation (A-R), capture all of them
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

# Snapshot the original values (and D-column date formatting) for every row before overwriting
$snapshot = @{}
foreach ($row in 2..22) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value()
    }
    $snapshot[$row] = $rowData
}

# Apply the permutation: destination row gets the data that used to live at source row
foreach ($destRow in 2..22) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcData[$col]
    }
}
